$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$helper = $ws.Cells.Item(100, 1)

$ws.Cells.Item(2, 5).Value = "2026-02-26 05:48:27"
$ws.Cells.Item(2, 13).Value = "2.9 °C 5:29 TU"
$ws.Cells.Item(3, 5).Value = "2026-02-26 05:48:30"
$ws.Cells.Item(3, 14).Value = "0.4 °C 5:01 TU"
$ws.Cells.Item(3, 15).Value = "1.2 °C"
$ws.Cells.Item(4, 5).Value = "2026-02-26 05:48:32"
$ws.Cells.Item(4, 10).Value = "1026.6 hPa"
$ws.Cells.Item(4, 15).Value = "6.6 °C"
$ws.Cells.Item(5, 5).Value = "2026-02-26 05:48:35"
$ws.Cells.Item(5, 14).Value = "2.6 °C 5:16 TU"
$ws.Cells.Item(5, 15).Value = "3.3 °C"
$ws.Cells.Item(6, 5).Value = "2026-02-26 05:48:37"
$ws.Cells.Item(7, 5).Value = "2026-02-26 05:48:40"
$helper.NumberFormat = "@"
$helper.Value = "89%"
$helper.Copy()
$ws.Cells.Item(7, 8).PasteSpecial(-4163)
$helper.Clear()
$ws.Cells.Item(7, 10).Value = "1026.0 hPa"
$ws.Cells.Item(7, 11).Value = "-0.1 MJ/m2"
$ws.Cells.Item(7, 14).Value = "10.9 °C 5:17 TU"
$ws.Cells.Item(8, 5).Value = "2026-02-26 05:48:42"
$helper.NumberFormat = "@"
$helper.Value = "95%"
$helper.Copy()
$ws.Cells.Item(8, 8).PasteSpecial(-4163)
$helper.Clear()
$ws.Cells.Item(8, 14).Value = "8.7 °C 5:00 TU"
$ws.Cells.Item(9, 5).Value = "2026-02-26 05:48:45"
$ws.Cells.Item(9, 14).Value = "9.2 °C 5:00 TU"
$ws.Cells.Item(9, 15).Value = "10.4 °C"
$ws.Cells.Item(10, 5).Value = "2026-02-26 05:48:48"
$ws.Cells.Item(10, 13).Value = "6.0 °C 5:29 TU"
$ws.Cells.Item(10, 15).Value = "4.5 °C"
$ws.Cells.Item(11, 5).Value = "2026-02-26 05:48:50"
$ws.Cells.Item(11, 15).Value = "1.9 °C"
$ws.Cells.Item(12, 5).Value = "2026-02-26 05:48:53"
$ws.Cells.Item(12, 14).Value = "8.1 °C 5:26 TU"
$ws.Cells.Item(12, 15).Value = "9.5 °C"
$ws.Cells.Item(13, 5).Value = "2026-02-26 05:48:55"
$helper.NumberFormat = "@"
$helper.Value = "94%"
$helper.Copy()
$ws.Cells.Item(13, 8).PasteSpecial(-4163)
$helper.Clear()
$ws.Cells.Item(13, 10).Value = "1032.1 hPa"
$ws.Cells.Item(13, 15).Value = "-1.2 °C"
$ws.Cells.Item(14, 5).Value = "2026-02-26 05:48:58"
$ws.Cells.Item(14, 14).Value = "8.5 °C 5:10 TU"
$ws.Cells.Item(14, 15).Value = "9.6 °C"
$ws.Cells.Item(15, 5).Value = "2026-02-26 05:48:59"
$helper.NumberFormat = "@"
$helper.Value = "95%"
$helper.Copy()
$ws.Cells.Item(15, 8).PasteSpecial(-4163)
$helper.Clear()
$ws.Cells.Item(15, 14).Value = "7.8 °C 5:29 TU"
$ws.Cells.Item(15, 15).Value = "10.0 °C"
$ws.Cells.Item(16, 5).Value = "2026-02-26 05:49:00"
$ws.Cells.Item(17, 5).Value = "2026-02-26 05:49:01"
$ws.Cells.Item(17, 14).Value = "5.5 °C 5:10 TU"
$ws.Cells.Item(18, 5).Value = "2026-02-26 05:49:02"
$ws.Cells.Item(18, 10).Value = "1026.5 hPa"
$ws.Cells.Item(18, 14).Value = "7.0 °C 5:04 TU"
$ws.Cells.Item(18, 15).Value = "8.2 °C"
$ws.Cells.Item(19, 5).Value = "2026-02-26 05:49:03"
$helper.NumberFormat = "@"
$helper.Value = "63%"
$helper.Copy()
$ws.Cells.Item(19, 8).PasteSpecial(-4163)
$helper.Clear()
$ws.Cells.Item(19, 12).Value = "7.6 km/h - 325º 5:10 TU"
$ws.Cells.Item(19, 15).Value = "7.7 °C"
$ws.Cells.Item(20, 5).Value = "2026-02-26 05:49:05"
$ws.Cells.Item(21, 5).Value = "2026-02-26 05:49:06"
$helper.NumberFormat = "@"
$helper.Value = "85%"
$helper.Copy()
$ws.Cells.Item(21, 8).PasteSpecial(-4163)
$helper.Clear()
$ws.Cells.Item(21, 10).Value = "1029.0 hPa"
$ws.Cells.Item(21, 14).Value = "2.0 °C 5:23 TU"
$ws.Cells.Item(21, 15).Value = "4.0 °C"
$ws.Cells.Item(22, 5).Value = "2026-02-26 05:49:08"
$helper.NumberFormat = "@"
$helper.Value = "56%"
$helper.Copy()
$ws.Cells.Item(22, 8).PasteSpecial(-4163)
$helper.Clear()
$ws.Cells.Item(23, 5).Value = "2026-02-26 05:49:10"
$ws.Cells.Item(23, 14).Value = "0.9 °C 5:01 TU"
$ws.Cells.Item(23, 15).Value = "2.2 °C"
$ws.Cells.Item(24, 5).Value = "2026-02-26 05:49:13"
$helper.NumberFormat = "@"
$helper.Value = "84%"
$helper.Copy()
$ws.Cells.Item(24, 8).PasteSpecial(-4163)
$helper.Clear()
$ws.Cells.Item(24, 10).Value = "1026.2 hPa"
$ws.Cells.Item(24, 14).Value = "1.1 °C 5:25 TU"
$ws.Cells.Item(24, 15).Value = "6.1 °C"
$ws.Cells.Item(25, 5).Value = "2026-02-26 05:49:15"
$helper.NumberFormat = "@"
$helper.Value = "39%"
$helper.Copy()
$ws.Cells.Item(25, 8).PasteSpecial(-4163)
$helper.Clear()
$ws.Cells.Item(25, 15).Value = "2.9 °C"
$ws.Cells.Item(26, 5).Value = "2026-02-26 05:49:18"
$helper.NumberFormat = "@"
$helper.Value = "42%"
$helper.Copy()
$ws.Cells.Item(26, 8).PasteSpecial(-4163)
$helper.Clear()
$ws.Cells.Item(26, 13).Value = "8.8 °C 5:27 TU"
$ws.Cells.Item(26, 15).Value = "7.4 °C"
$ws.Cells.Item(27, 5).Value = "2026-02-26 05:49:20"
$helper.NumberFormat = "@"
$helper.Value = "55%"
$helper.Copy()
$ws.Cells.Item(27, 8).PasteSpecial(-4163)
$helper.Clear()
$ws.Cells.Item(27, 14).Value = "1.5 °C 5:18 TU"
$ws.Cells.Item(28, 5).Value = "2026-02-26 05:49:23"
$ws.Cells.Item(28, 10).Value = "1026.3 hPa"
$ws.Cells.Item(28, 14).Value = "6.9 °C 5:29 TU"
$ws.Cells.Item(28, 15).Value = "8.2 °C"
$ws.Cells.Item(29, 5).Value = "2026-02-26 05:49:25"
$helper.NumberFormat = "@"
$helper.Value = "97%"
$helper.Copy()
$ws.Cells.Item(29, 8).PasteSpecial(-4163)
$helper.Clear()
$ws.Cells.Item(29, 12).Value = "13.7 km/h - 350º 5:01 TU"
$ws.Cells.Item(30, 5).Value = "2026-02-26 05:49:28"
$ws.Cells.Item(30, 14).Value = "9.6 °C 5:17 TU"
$ws.Cells.Item(30, 15).Value = "10.6 °C"
$ws.Cells.Item(31, 5).Value = "2026-02-26 05:49:30"
$ws.Cells.Item(31, 10).Value = "1025.8 hPa"
$ws.Cells.Item(32, 5).Value = "2026-02-26 05:49:33"
$helper.NumberFormat = "@"
$helper.Value = "77%"
$helper.Copy()
$ws.Cells.Item(32, 8).PasteSpecial(-4163)
$helper.Clear()
$ws.Cells.Item(32, 14).Value = "-1.3 °C 5:26 TU"
$ws.Cells.Item(32, 15).Value = "0.9 °C"
$ws.Cells.Item(33, 5).Value = "2026-02-26 05:49:35"
$ws.Cells.Item(33, 10).Value = "1029.6 hPa"
$ws.Cells.Item(33, 14).Value = "0.7 °C 5:09 TU"
$ws.Cells.Item(33, 15).Value = "2.3 °C"
$ws.Cells.Item(34, 5).Value = "2026-02-26 05:49:38"
$helper.NumberFormat = "@"
$helper.Value = "47%"
$helper.Copy()
$ws.Cells.Item(34, 8).PasteSpecial(-4163)
$helper.Clear()
$ws.Cells.Item(34, 15).Value = "3.1 °C"
$ws.Cells.Item(35, 5).Value = "2026-02-26 05:49:41"
$ws.Cells.Item(35, 10).Value = "1025.5 hPa"
$ws.Cells.Item(35, 14).Value = "7.2 °C 5:07 TU"
$ws.Cells.Item(35, 15).Value = "8.7 °C"
$ws.Cells.Item(36, 5).Value = "2026-02-26 05:49:43"
$ws.Cells.Item(37, 5).Value = "2026-02-26 05:49:46"
$ws.Cells.Item(37, 10).Value = "1029.1 hPa"
$ws.Cells.Item(37, 14).Value = "0.8 °C 5:29 TU"
$ws.Cells.Item(37, 15).Value = "2.1 °C"
$ws.Cells.Item(38, 5).Value = "2026-02-26 05:49:48"
$ws.Cells.Item(39, 5).Value = "2026-02-26 05:49:50"
$ws.Cells.Item(39, 15).Value = "2.7 °C"
$ws.Cells.Item(40, 5).Value = "2026-02-26 05:49:53"
$ws.Cells.Item(40, 10).Value = "1029.7 hPa"
$ws.Cells.Item(40, 15).Value = "2.1 °C"
$ws.Cells.Item(41, 5).Value = "2026-02-26 05:49:56"
$helper.NumberFormat = "@"
$helper.Value = "97%"
$helper.Copy()
$ws.Cells.Item(41, 8).PasteSpecial(-4163)
$helper.Clear()
$ws.Cells.Item(41, 14).Value = "5.4 °C 5:29 TU"
$ws.Cells.Item(41, 15).Value = "7.8 °C"
$ws.Cells.Item(42, 5).Value = "2026-02-26 05:49:58"
$ws.Cells.Item(43, 5).Value = "2026-02-26 05:50:01"
$ws.Cells.Item(43, 11).Value = "-0.1 MJ/m2"
$ws.Cells.Item(43, 14).Value = "1.1 °C 5:05 TU"
$ws.Cells.Item(43, 15).Value = "2.7 °C"
$ws.Cells.Item(44, 5).Value = "2026-02-26 05:50:03"
$ws.Cells.Item(44, 15).Value = "0.2 °C"
$ws.Cells.Item(45, 5).Value = "2026-02-26 05:50:06"
$ws.Cells.Item(45, 10).Value = "1027.5 hPa"
$ws.Cells.Item(45, 14).Value = "4.4 °C 5:13 TU"
$ws.Cells.Item(45, 15).Value = "6.1 °C"
$ws.Cells.Item(46, 5).Value = "2026-02-26 05:50:08"
$ws.Cells.Item(46, 10).Value = "1026.1 hPa"
$ws.Cells.Item(46, 14).Value = "4.1 °C 5:02 TU"
$ws.Cells.Item(46, 15).Value = "7.1 °C"
$excel.CutCopyMode = 0
